$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new timesheet entry as row 17: same date as row 16, 2.5 hours,
# new comment "DefaultSerializer + testy" (added to shared strings).
$ws.Range("A17").Value = 44025
$ws.Range("B17").Value = 2.5
$ws.Range("C17").Value = "DefaultSerializer + testy"

# Move the selection to the newly added comment cell, as in the saved file.
$ws.Range("C17").Select()
